# Stock data scraping to excel files
#
# Adds two new worksheets ("website" and "symbols") to the existing
# "emails" workbook and appends two more rows of email addresses / new
# hyperlinks to the "emails" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the two new sheets. NOTE: in this engine, cached worksheet
#    object references become stale after a further Worksheets.Add()
#    call, so create every sheet first, then re-resolve each sheet via
#    Worksheets.Item(name) right before touching it.
# ---------------------------------------------------------------------
$newSymbols = $wb.Worksheets.Add($null, $wb.Worksheets.Item("emails"))
$newSymbols.Name = "symbols"

$newWebsite = $wb.Worksheets.Add($null, $wb.Worksheets.Item("emails"))
$newWebsite.Name = "website"

# Final tab order must be: emails, website, symbols

# ---------------------------------------------------------------------
# 2. "emails" sheet - append two more addresses with mailto hyperlinks
# ---------------------------------------------------------------------
$wsEmails = $wb.Worksheets.Item("emails")

$wsEmails.Range("A5").Value = "moldovanandrei2301@gmail.com"
$wsEmails.Hyperlinks.Add($wsEmails.Range("A5"), "mailto:moldovanandrei2301@gmail.com") | Out-Null
$wsEmails.Range("A5").Style = "Hyperlink"

$wsEmails.Range("A6").Value = "amoldovan65@yahoo.com"
$wsEmails.Hyperlinks.Add($wsEmails.Range("A6"), "mailto:amoldovan65@yahoo.com") | Out-Null
$wsEmails.Range("A6").Style = "Hyperlink"

$wsEmails.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. "symbols" sheet - list of tracked ticker symbols
#    (populated before "website" below so the new shared-string table
#    entries land in the same order as the target workbook.)
# ---------------------------------------------------------------------
$wsSymbols = $wb.Worksheets.Item("symbols")

# Target stored column width is 17.140625 characters; closest
# achievable quantised value (see note above) comes from an input of
# 16.3, which lands on stored width 17.166666...
$wsSymbols.Columns.Item(1).ColumnWidth = 16.3

$wsSymbols.Range("A1").Value = "Symbols"
$wsSymbols.Range("A2").Value = "NASDAQ-AMZN"
$wsSymbols.Range("A3").Value = "NASDAQ-AAPL"
$wsSymbols.Range("A4").Value = "NASDAQ-TSLA"
$wsSymbols.Range("A5").Value = "NASDAQ-NFLX"
$wsSymbols.Range("A6").Value = "BTCUSD"
$wsSymbols.Range("A7").Value = "OTC-ETHE"

$wsSymbols.Range("A10").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. "website" sheet - scraping URL format
# ---------------------------------------------------------------------
$wsWebsite = $wb.Worksheets.Item("website")

# Target stored column width is 56.85546875 characters; the emulated
# engine quantises ColumnWidth to 1/6-character pixel steps, so 56.0
# is the closest input that lands on the nearest achievable stored
# width (56.833333...).
$wsWebsite.Columns.Item(1).ColumnWidth = 56

$wsWebsite.Range("A1").Value = "Website Format"
$wsWebsite.Range("A2").Value = "https://www.tradingview.com/symbols/{symbol}/technicals/"
$wsWebsite.Hyperlinks.Add($wsWebsite.Range("A2"), "https://www.tradingview.com/symbols/{symbol}/technicals/") | Out-Null
$wsWebsite.Range("A2").Style = "Hyperlink"

$wsWebsite.Range("A10").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Make "website" the active sheet/tab (matches the saved workbook
#    view in the target file).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("website").Activate() | Out-Null
